$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update existing "Sheet4" tab (A1:M23 table) with two newly-measured rows
#    of data (B21/C21 and B22/C22), which feed the existing K/L "delta" shared
#    formulas for rows 21 and 22.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Activate()

$ws4.Range("B21").Value = -17.8
$ws4.Range("B21").Font.Name = "Menlo"
$ws4.Range("B21").Font.Size = 11
$ws4.Range("B21").Font.Color = 0

$ws4.Range("C21").Value = -48
$ws4.Range("C21").Font.Name = "Menlo"
$ws4.Range("C21").Font.Size = 11
$ws4.Range("C21").Font.Color = 0

$ws4.Range("B22").Value = -18.1
$ws4.Range("B22").Font.Name = "Menlo"
$ws4.Range("B22").Font.Size = 11
$ws4.Range("B22").Font.Color = 0

$ws4.Range("C22").Value = -37.7
$ws4.Range("C22").Font.Name = "Menlo"
$ws4.Range("C22").Font.Size = 11
$ws4.Range("C22").Font.Color = 0

$null = $ws4.Range("K9").Select()

# ---------------------------------------------------------------------------
# 2. Add a brand new worksheet "Sheet6" at the end of the workbook containing
#    a fresh vision-data "fudge factor" table for a clean field.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws6.Name = "Sheet6"

# Header row (all of these reuse pre-existing shared strings, except E1 which
# is deferred below so the new shared strings come out in the same order as
# the original authoring session).
$ws6.Range("A1").Value = "COORD"
$ws6.Range("B1").Value = "Predicted X"
$ws6.Range("C1").Value = "Predicted Y"
$ws6.Range("F1").Value = "Real Y"
$ws6.Range("K1").Value = "Delta X"
$ws6.Range("L1").Value = "Delta Y"
$ws6.Range("M1").Value = "Delta H (deg)"

function Style-Menlo($rng) {
    $rng.Font.Name = "Menlo"
    $rng.Font.Size = 11
    $rng.Font.Color = 0
}

# Row 2
$ws6.Range("A2").Value = "B(2, 3)"
$ws6.Range("B2").Value = 27
Style-Menlo $ws6.Range("B2")
$ws6.Range("C2").Value = -38.6
Style-Menlo $ws6.Range("C2")
Style-Menlo $ws6.Range("D2")
$ws6.Range("E2").Value = 28
$ws6.Range("F2").Value = 39
$ws6.Range("K2").Formula = "= E2 - ABS(B2)"
$ws6.Range("L2").Formula = "= F2 - ABS(C2)"
$ws6.Range("M2").Formula = "= G2 - ABS(D2)"

# Row 3
$ws6.Range("A3").Value = "B(2, 2)"
$ws6.Range("B3").Value = 34.6
Style-Menlo $ws6.Range("B3")
$ws6.Range("C3").Value = -38.5
Style-Menlo $ws6.Range("C3")
$ws6.Range("E3").Value = 35.5
$ws6.Range("F3").Value = 39
$ws6.Range("K3").Formula = "= E3 - ABS(B3)"
$ws6.Range("L3").Formula = "= F3 - ABS(C3)"
$ws6.Range("M3").Formula = "= G3 - ABS(D3)"

# Row 4
$ws6.Range("A4").Value = "B(2, 1)"
$ws6.Range("B4").Value = 41.7
Style-Menlo $ws6.Range("B4")
$ws6.Range("C4").Value = -37.8
Style-Menlo $ws6.Range("C4")
$ws6.Range("E4").Value = 43
$ws6.Range("F4").Value = 38.5
$ws6.Range("K4").Formula = "= E4 - ABS(B4)"
$ws6.Range("L4").Formula = "= F4 - ABS(C4)"
$ws6.Range("M4").Formula = "= G4 - ABS(D4)"

# Row 7
$ws6.Range("A7").Value = "B(0, 3)"
$ws6.Range("B7").Value = 0.1
Style-Menlo $ws6.Range("B7")
$ws6.Range("C7").Value = 8.7
Style-Menlo $ws6.Range("C7")
$ws6.Range("E1").Value = "Real X (TAG REL)"
$ws6.Range("E7").Value = 0
$ws6.Range("F7").Value = 8
$ws6.Range("K7").Formula = "= E7 - ABS(B7)"
$ws6.Range("L7").Formula = "= F7 - ABS(C7)"
$ws6.Range("M7").Formula = "= G7 - ABS(D7)"

# Row 8
$ws6.Range("A8").Value = "B(1, 3)"
$ws6.Range("B8").Value = 0.1
Style-Menlo $ws6.Range("B8")
$ws6.Range("C8").Value = 14.8
Style-Menlo $ws6.Range("C8")
$ws6.Range("E8").Value = 0
$ws6.Range("F8").Value = 14
$ws6.Range("K8").Formula = "= E8 - ABS(B8)"
$ws6.Range("L8").Formula = "= F8 - ABS(C8)"
$ws6.Range("M8").Formula = "= G8 - ABS(D8)"

# Row 9
$ws6.Range("A9").Value = "B(2, 3)"
$ws6.Range("B9").Value = 0
Style-Menlo $ws6.Range("B9")
$ws6.Range("C9").Value = 24.9
Style-Menlo $ws6.Range("C9")
$ws6.Range("E9").Value = 0
$ws6.Range("F9").Value = 24
$ws6.Range("K9").Formula = "= E9 - ABS(B9)"
$ws6.Range("L9").Formula = "= F9 - ABS(C9)"
$ws6.Range("M9").Formula = "= G9 - ABS(D9)"

# Row 10
$ws6.Range("A10").Value = "B(0, 2)"
$ws6.Range("B10").Value = 0.1
Style-Menlo $ws6.Range("B10")
$ws6.Range("C10").Value = 8.8
Style-Menlo $ws6.Range("C10")
$ws6.Range("E10").Value = 0
$ws6.Range("F10").Value = 8.25
$ws6.Range("K10").Formula = "= E10 - ABS(B10)"
$ws6.Range("L10").Formula = "= F10 - ABS(C10)"
$ws6.Range("M10").Formula = "= G10 - ABS(D10)"

# Row 11
$ws6.Range("A11").Value = "B(1, 2)"
$ws6.Range("B11").Value = -0.3
Style-Menlo $ws6.Range("B11")
$ws6.Range("C11").Value = 14.9
Style-Menlo $ws6.Range("C11")
$ws6.Range("E11").Value = 0
$ws6.Range("F11").Value = 14.25
$ws6.Range("K11").Formula = "= E11 - ABS(B11)"
$ws6.Range("L11").Formula = "= F11 - ABS(C11)"

# Row 12
$ws6.Range("A12").Value = "B(2, 2)"
$ws6.Range("B12").Value = -0.6
Style-Menlo $ws6.Range("B12")
$ws6.Range("C12").Value = 25
Style-Menlo $ws6.Range("C12")
$ws6.Range("E12").Value = 0
$ws6.Range("F12").Value = 24.5
$ws6.Range("K12").Formula = "= E12 - ABS(B12)"
$ws6.Range("L12").Formula = "= F12 - ABS(C12)"

# Row 13
$ws6.Range("A13").Value = "B(0, 1)"
$ws6.Range("B13").Value = 0
Style-Menlo $ws6.Range("B13")
$ws6.Range("C13").Value = 9
Style-Menlo $ws6.Range("C13")
$ws6.Range("E13").Value = 0
$ws6.Range("F13").Value = 8
$ws6.Range("K13").Formula = "= E13 - ABS(B13)"
$ws6.Range("L13").Formula = "= F13 - ABS(C13)"

# Row 14
$ws6.Range("A14").Value = "B(1, 1)"
$ws6.Range("B14").Value = 0.1
Style-Menlo $ws6.Range("B14")
$ws6.Range("C14").Value = 15
Style-Menlo $ws6.Range("C14")
$ws6.Range("E14").Value = 0
$ws6.Range("F14").Value = 14

# Row 16
$ws6.Range("L16").Value = "Fudge Factor:"

# Row 20
$ws6.Range("L20").Value = 1
$ws6.Range("M20").Formula = "=AVERAGE(L2:L4, L7:L13)"

# Rows 22-32 (fudge-corrected deltas)
$ws6.Range("L22").Formula = "=L2-`$M`$20"
$ws6.Range("L23").Formula = "=L3-`$M`$20"
$ws6.Range("L24").Formula = "=L4-`$M`$20"
$ws6.Range("L27").Formula = "=L7-`$M`$20"
$ws6.Range("L28").Formula = "=L8-`$M`$20"
$ws6.Range("L29").Formula = "=L9-`$M`$20"
$ws6.Range("L30").Formula = "=L10-`$M`$20"
$ws6.Range("L31").Formula = "=L11-`$M`$20"
$ws6.Range("L32").Formula = "=L12-`$M`$20"

$ws6.Activate()
$null = $ws6.Range("N13").Select()
